# Refresh the cryptos price list: updated Price (D) / hourly-Volume% (E) figures, and
# two rank swaps where Coin (B) + Link (C) also change:
#   rows 35/36  HuobiToken <-> Frax
#   rows 40/41  TheSandbox <-> Algorand
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume cells are stored as literal text (e.g. "0.3610", "1.001",
# "  -4.05%  "). Values that are valid numbers get a leading apostrophe (the
# Excel text-qualifier) so they are kept as typed instead of being auto-converted
# to a Number (which would also silently drop significant trailing zeros, e.g.
# "0.3610" -> 0.361). The apostrophe itself is not part of the stored value.

$ws.Range('D2').Value = '26.319.95'
$ws.Range('D3').Value = '1.758.73'
$ws.Range('E3').Value = '  -4.05%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D6').Value = "'306.04"
$ws.Range('E6').Value = '  -1.88%  '
$ws.Range('D7').Value = "'0.4282"
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').Value = "'0.3610"
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('D9').Value = "'0.07025"
$ws.Range('E9').Value = '  -3.17%  '
$ws.Range('D10').Value = "'0.8301"
$ws.Range('E10').Value = '  -3.75%  '
$ws.Range('D11').Value = "'20.09"
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('D12').Value = '1.730.18'
$ws.Range('E12').Value = '  -4.10%  '
$ws.Range('D13').Value = "'5.207"
$ws.Range('E13').Value = '  -3.72%  '
$ws.Range('D14').Value = "'6.345"
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').Value = "'0.06785"
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').Value = "'1.002"
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = "'78.99"
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').Value = "'0.000008642"
$ws.Range('E18').Value = '  -2.77%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = "'14.84"
$ws.Range('E20').Value = '  -3.52%  '
$ws.Range('D21').Value = '26.062.32'
$ws.Range('E21').Value = '  -5.04%  '
$ws.Range('D22').Value = "'4.974"
$ws.Range('E22').Value = '  -3.34%  '
$ws.Range('E23').Value = '  +1.91%  '
$ws.Range('D24').Value = '1.957.17'
$ws.Range('E24').Value = '  -4.64%  '
$ws.Range('E25').Value = '  -4.54%  '
$ws.Range('D26').Value = "'151.81"
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').Value = "'18.10"
$ws.Range('E27').Value = '  -3.71%  '
$ws.Range('D28').Value = "'114.43"
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = "'5.015"
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('D30').Value = "'1.675"
$ws.Range('E30').Value = '  -7.88%  '
$ws.Range('D31').Value = "'0.08880"
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').Value = "'0.7184"
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('D33').Value = "'4.295"
$ws.Range('E33').Value = '  -5.20%  '
$ws.Range('D34').Value = "'1.097"
$ws.Range('E34').Value = '  -2.74%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = "'1.001"
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.740"
$ws.Range('E36').Value = '  -8.17%  '
$ws.Range('D37').Value = "'1.064"
$ws.Range('E37').Value = '  -2.33%  '
$ws.Range('D38').Value = "'0.05079"
$ws.Range('E38').Value = '  -4.19%  '
$ws.Range('D39').Value = "'0.01874"
$ws.Range('E39').Value = '  -3.05%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = "'0.1596"
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'0.4873"
$ws.Range('E41').Value = '  -3.90%  '
$ws.Range('D42').Value = "'2.474"
$ws.Range('E42').Value = '  -11.56%  '
$ws.Range('D43').Value = "'6.136"
$ws.Range('E43').Value = '  -5.25%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = "'7.952"
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('D46').Value = "'104.51"
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('D47').Value = "'1.001"
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = "'9.964"
$ws.Range('E48').Value = '  -4.42%  '
$ws.Range('E49').Value = '  -4.99%  '
$ws.Range('D50').Value = "'0.4453"
$ws.Range('E50').Value = '  -4.51%  '
$ws.Range('D51').Value = "'1.564"
$ws.Range('E51').Value = '  -3.14%  '
